$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.537.72"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "'2.107.39"
$ws.Range("E3").Value = "  +4.70%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'329.24"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.5261"
$ws.Range("E7").Value = "  +2.60%  "
$ws.Range("D8").Value = "'0.4376"
$ws.Range("E8").Value = "  +2.81%  "
$ws.Range("D9").Value = "'0.08884"
$ws.Range("D10").Value = "'47.54"
$ws.Range("E10").Value = "  +9.73%  "
$ws.Range("D11").Value = "'1.164"
$ws.Range("E11").Value = "  +2.53%  "
$ws.Range("D12").Value = "'24.62"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").Value = "'2.107.08"
$ws.Range("E13").Value = "  +4.71%  "
$ws.Range("D14").Value = "'6.734"
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("E15").Value = "  +4.02%  "
$ws.Range("D16").Value = "'96.31"
$ws.Range("E16").Value = "  +2.16%  "
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("D19").Value = "'0.06637"
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("D20").Value = "'19.03"
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "'6.317"
$ws.Range("E22").Value = "  +1.78%  "
$ws.Range("D23").Value = "'30.594.83"
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("D24").Value = "'12.28"
$ws.Range("E24").Value = "  +3.97%  "
$ws.Range("D25").Value = "'2.357"
$ws.Range("E25").Value = "  +3.95%  "
$ws.Range("D26").Value = "'2.360.35"
$ws.Range("E26").Value = "  +4.97%  "
$ws.Range("D27").Value = "'22.43"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +7.35%  "
$ws.Range("D29").Value = "'161.97"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").Value = "'132.73"
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("D31").Value = "'1.206"
$ws.Range("E31").Value = "  +5.76%  "
$ws.Range("E32").Value = "  +2.32%  "
$ws.Range("D33").Value = "'1.682"
$ws.Range("E33").Value = "  +22.43%  "
$ws.Range("D34").Value = "'6.196"
$ws.Range("E34").Value = "  +1.87%  "
$ws.Range("E35").Value = "  +2.63%  "
$ws.Range("D36").Value = "'10.19"
$ws.Range("E36").Value = "  +11.41%  "
$ws.Range("D37").Value = "'0.02582"
$ws.Range("E37").Value = "  +2.29%  "
$ws.Range("E38").Value = "  +0.72%  "
$ws.Range("D39").Value = "'12.73"
$ws.Range("E39").Value = "  +2.81%  "
$ws.Range("D40").Value = "'0.06684"
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").Value = "'0.2289"
$ws.Range("E41").Value = "  +4.27%  "
$ws.Range("D42").Value = "'0.6844"
$ws.Range("E42").Value = "  +2.79%  "
$ws.Range("D43").Value = "'1.259"
$ws.Range("E43").Value = "  +1.46%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "'14.07"
$ws.Range("E45").Value = "  +2.57%  "
$ws.Range("D46").Value = "'0.6379"
$ws.Range("E46").Value = "  +3.57%  "
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("D48").Value = "'3.624"
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("D49").Value = "'1.253"
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").Value = "'1.201"
$ws.Range("E50").Value = "  +8.44%  "
$ws.Range("D51").Value = "'82.41"
$ws.Range("E51").Value = "  +2.18%  "
